$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.103.72"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "'2.638.22"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'523.77"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("D6").Value = "'146.18"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'2.656.42"
$ws.Range("E9").Value = "  +2.71%  "

$ws.Range("D10").Value = "'6.33"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("D14").Value = "'3.101.13"
$ws.Range("E14").Value = "  +2.61%  "

$ws.Range("D15").Value = "'59.084.19"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").Value = "'20.99"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'2.641.50"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").Value = "'347.18"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("E20").Value = "  -1.37%  "

$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("D22").Value = "'6.18"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Value = "'61.89"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("D27").Value = "'0.993"
$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").Value = "'7.12"
$ws.Range("E29").Value = "  +0.89%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").Value = "'6.27"
$ws.Range("E31").Value = "  +4.50%  "

$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("E33").Value = "  +1.03%  "

$ws.Range("D34").Value = "'150.55"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").Value = "'0.983"
$ws.Range("E35").Value = "  +5.50%  "

$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "'36.80"

$ws.Range("D39").Value = "'0.851"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("E40").Value = "  +2.45%  "

$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("D42").Value = "'279.00"
$ws.Range("E42").Value = "  -4.20%  "

$ws.Range("D43").Value = "'0.612"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0988"
$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.995"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").Value = "'19.61"
$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("D47").Value = "'0.0524"
$ws.Range("E47").Value = "  -3.12%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.75"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("D49").Value = "'10.30"
$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0230"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.996.49"
$ws.Range("E51").Value = "  +3.32%  "
